$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DI34")
$fmtSrc = $wb.Worksheets.Item("Master_DB").Range("C2")

# Helper: apply the existing date-format style (xf index 3, numFmtId 14)
# to a cell by copying formats only from a donor cell that already uses
# it, so no new cell style gets minted.
function Set-DateStyle($range) {
    $fmtSrc.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

# Helper: write literal text into a date-formatted cell without Excel
# re-interpreting it as a date serial. We evaluate it as a formula that
# yields a text string, then freeze it back down to a plain value via a
# values-only paste (keeps the already-applied number format/style).
function Set-LiteralText($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# --- New text values are entered in this order so the workbook's shared-
# string table grows in the same sequence as the authored edit. ---

# Header row: new "Description" column before the last column.
$ws.Range("I1").Value = "Description"

# Row 14 memo.
$ws.Range("J14").Value = "Deposit Cheque 212"

# Row 13 memo.
$ws.Range("J13").Value = "Deposit Cheque 111"

# Rows 15/16 share a literal (non-date) "05-20-2025" posting-date label.
Set-DateStyle $ws.Range("B15")
Set-LiteralText $ws.Range("B15") "05-20-2025"
Set-DateStyle $ws.Range("B16")
Set-LiteralText $ws.Range("B16") "05-20-2025"

# Row 15 memo.
$ws.Range("J15").Value = "Deposit Cheque 313"

# Row 12: add a "Charges" note in column J.
$ws.Range("J12").Value = "Charges"

# --- Remaining cell values for the four new/changed rows ---

# Row 13 (SSN 222222222 stays).
$ws.Range("A13").Value = 222222222
Set-DateStyle $ws.Range("B13")
$ws.Range("B13").Value = 45996
$ws.Range("C13").Value = 1003321229
$ws.Range("D13").Value = "SVG"
Set-DateStyle $ws.Range("E13")
$ws.Range("E13").Value = 45801
$ws.Range("F13").Value = 100
$ws.Range("G13").Value = "C"
$ws.Range("H13").Value = 950

# Row 14 (SSN changes 333333333 -> 222222222).
$ws.Range("A14").Value = 222222222
Set-DateStyle $ws.Range("B14")
$ws.Range("B14").Value = 45996
$ws.Range("C14").Value = 1003321229
$ws.Range("D14").Value = "SVG"
Set-DateStyle $ws.Range("E14")
$ws.Range("E14").Value = 45801
$ws.Range("F14").Value = 50
$ws.Range("G14").Value = "C"
$ws.Range("H14").Value = 1000

# Row 15 (SSN changes 444444444 -> 222222222).
$ws.Range("A15").Value = 222222222
$ws.Range("C15").Value = 1003321229
$ws.Range("D15").Value = "SVG"
Set-DateStyle $ws.Range("E15")
$ws.Range("E15").Value = 45801
$ws.Range("F15").Value = 100
$ws.Range("G15").Value = "C"
$ws.Range("H15").Value = 1100

# Row 16 (SSN changes 666666666 -> 222222222).
$ws.Range("A16").Value = 222222222
$ws.Range("C16").Value = 1003321229
$ws.Range("D16").Value = "SVG"
Set-DateStyle $ws.Range("E16")
$ws.Range("E16").Value = 45801
$ws.Range("F16").Value = 1100
$ws.Range("G16").Value = "D"
$ws.Range("H16").Value = 0

# --- Selection moves to the newly added block ---
$ws.Range("A17:A19").Select()
